$d = $word.ActiveDocument

function Replace-WithBold {
    param(
        [string]$SearchText,
        [string]$NewText,
        [string]$BoldText
    )

    $r = $d.Content
    $found = $r.Find.Execute($SearchText, $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $SearchText
        return
    }

    $start = $r.Start
    $end = $r.End

    $full = $d.Range($start, $end)
    $full.Text = $NewText

    if ($BoldText) {
        $boldIdx = $NewText.IndexOf($BoldText)
        if ($boldIdx -ge 0) {
            $boldStart = $start + $boldIdx
            $boldEnd = $boldStart + $BoldText.Length
            $boldRange = $d.Range($boldStart, $boldEnd)
            $boldRange.Bold = 1
            $boldRange.BoldBi = 1
        }
    }
}

# "Play on words" row
Replace-WithBold "No play on words" "The tweet has a No play on words" "No"
Replace-WithBold "Minor play on words" "The tweet has a Minor play on words." "Minor"
Replace-WithBold "Cleaver play on words" "The tweet has a Cleaver play on words." "Cleaver"

# "Level of Whitty(ness)" header fix -> "Level of Wittiness"
Replace-WithBold "Whitty" "Wittiness" ""

# Wittiness row
Replace-WithBold "No wittiness" "The tweet has No wittiness." "No"
Replace-WithBold "Low level of wittiness" "The tweet has a Low level of wittiness." "Low level"
Replace-WithBold "A good level of wittiness" "The tweet has a good level of wittiness." "good level"

# Laughter row
Replace-WithBold "No laughter" "The tweet gives you No laughter." "No"
Replace-WithBold "Some laughter" "The tweet gave you Some laughter." "Some"
Replace-WithBold "A lot of laughter" "The tweet gave you a Lot of laughter." "ot"

# Topic relevance row
Replace-WithBold "No relevant topic" "The tweet has No relevant topic." "No"
Replace-WithBold "Some topic relevance" "The tweet has Some topic relevance." "Some"
Replace-WithBold "A large level of topic relevance" "The tweet has a large level of topic relevance." "large level"

Write-Host "Rubric text updates complete"
